# Shrink the crossword grid from 7x7 (A1:G7) down to 5x5 (A1:E5):
# the two right-most columns (F:G) and the two bottom-most rows (6:7)
# are no longer part of the puzzle, so their contents are cleared
# (not structurally deleted -- the column-width runs in <cols> stay
# anchored at their original boundaries, only the used range shrinks).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("F1:G7").Clear()
$ws.Range("A6:G7").Clear()

# The former black ("blocked") squares inside the kept 5x5 area are no
# longer blocked -- turn them into regular "?" squares like the rest
# of the grid (clear their black-fill style, then give them the same
# shared-string value used everywhere else).
$ws.Range("A1:B1").ClearFormats()
$ws.Range("A2").ClearFormats()

$ws.Range("A1").Value = "?"
$ws.Range("B1").Value = "?"
$ws.Range("A2").Value = "?"

# Leave the cursor where the author left it after the edit.
$ws.Range("E2").Select() | Out-Null
